$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new To-do items (rows 32-37) first, so the new shared strings are
# appended to the shared string table in the same order the author typed
# them in.
$ws.Range("A32").Value = "Sort timestamp issue"
$ws.Range("A33").Value = "Fix buttons on random page"
$ws.Range("A34").Value = "Put scrollbar on divs where needed "
$ws.Range("A35").Value = "Email Andy and Ortho"
$ws.Range("A37").Value = "Admin adding events ? (possibly)"
$ws.Range("A36").Value = "Deletes - Delete messages etc ? (possibly)"

# Mark existing rows as "Done" in column B
$ws.Range("B3").Value = "Done"
$ws.Range("B4").Value = "Done"
$ws.Range("B5").Value = "Done"
$ws.Range("B6").Value = "Done"
$ws.Range("B7").Value = "Done"
$ws.Range("B17").Value = "Done"
$ws.Range("B32").Value = "Done"

# Mark existing rows as "In progress" in column B
$ws.Range("B29").Value = "In progress"
$ws.Range("B30").Value = "In progress"

# Update the active selection to match the final workbook state
$ws.Range("B16").Select()
